$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'" + '24.136.86'
$ws.Range("E2").Value = "'" + '  +10.44%  '
$ws.Range("E3").Value = "'" + '  +6.49%  '
$ws.Range("D4").Value = "'" + '1.006'
$ws.Range("E4").Value = "'" + '  +0.45%  '
$ws.Range("D5").Value = "'" + '308.16'
$ws.Range("E5").Value = "'" + '  +7.46%  '
$ws.Range("D6").Value = "'" + '1.000'
$ws.Range("E6").Value = "'" + '  +1.24%  '
$ws.Range("D7").Value = "'" + '0.3732'
$ws.Range("E7").Value = "'" + '  +0.99%  '
$ws.Range("D8").Value = "'" + '0.3441'
$ws.Range("E8").Value = "'" + '  +4.82%  '
$ws.Range("D9").Value = "'" + '47.43'
$ws.Range("E9").Value = "'" + '  +13.93%  '
$ws.Range("D10").Value = "'" + '1.188'
$ws.Range("E10").Value = "'" + '  +4.02%  '
$ws.Range("D11").Value = "'" + '0.07312'
$ws.Range("E11").Value = "'" + '  +3.85%  '
$ws.Range("D12").Value = "'" + '1.004'
$ws.Range("E12").Value = "'" + '  +0.51%  '
$ws.Range("D13").Value = "'" + '20.38'
$ws.Range("E13").Value = "'" + '  +1.54%  '
$ws.Range("D14").Value = "'" + '6.115'
$ws.Range("E14").Value = "'" + '  +4.56%  '
$ws.Range("D15").Value = "'" + '6.776'
$ws.Range("E15").Value = "'" + '  +3.70%  '
$ws.Range("D16").Value = "'" + '1.678.39'
$ws.Range("E16").Value = "'" + '  +6.64%  '
$ws.Range("E17").Value = "'" + '  +3.28%  '
$ws.Range("D18").Value = "'" + '1.001'
$ws.Range("E18").Value = "'" + '  +1.27%  '
$ws.Range("D19").Value = "'" + '0.06720'
$ws.Range("E19").Value = "'" + '  +6.13%  '
$ws.Range("D20").Value = "'" + '81.75'
$ws.Range("E20").Value = "'" + '  +8.38%  '
$ws.Range("E21").Value = "'" + '  +2.00%  '
$ws.Range("D22").Value = "'" + '6.115'
$ws.Range("E22").Value = "'" + '  +4.17%  '
$ws.Range("D23").Value = "'" + '12.02'
$ws.Range("E23").Value = "'" + '  +3.38%  '
$ws.Range("D24").Value = "'" + '24.122.59'
$ws.Range("E24").Value = "'" + '  +10.23%  '
$ws.Range("D25").Value = "'" + '2.417'
$ws.Range("E25").Value = "'" + '  +2.61%  '
$ws.Range("D26").Value = "'" + '3.366'
$ws.Range("E26").Value = "'" + '  -9.32%  '
$ws.Range("D27").Value = "'" + '2.661'
$ws.Range("D28").Value = "'" + '153.14'
$ws.Range("E28").Value = "'" + '  +2.15%  '
$ws.Range("D29").Value = "'" + '19.64'
$ws.Range("E29").Value = "'" + '  +5.93%  '
$ws.Range("D30").Value = "'" + '1.859.46'
$ws.Range("E30").Value = "'" + '  +6.40%  '
$ws.Range("D31").Value = "'" + '127.36'
$ws.Range("E31").Value = "'" + '  +5.85%  '
$ws.Range("D32").Value = "'" + '6.385'
$ws.Range("E32").Value = "'" + '  +16.85%  '
$ws.Range("D33").Value = "'" + '4.066'
$ws.Range("E33").Value = "'" + '  -1.83%  '
$ws.Range("D34").Value = "'" + '0.9772'
$ws.Range("E34").Value = "'" + '  +6.36%  '
$ws.Range("D35").Value = "'" + '1.761'
$ws.Range("E35").Value = "'" + '  +9.08%  '
$ws.Range("D36").Value = "'" + '0.08482'
$ws.Range("E36").Value = "'" + '  +3.22%  '
$ws.Range("D37").Value = "'" + '12.34'
$ws.Range("E37").Value = "'" + '  +5.59%  '
$ws.Range("D38").Value = "'" + '0.06434'
$ws.Range("E38").Value = "'" + '  +4.85%  '
$ws.Range("B39").Value = 'InternetComputer(DFINITY)'
$ws.Range("C39").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D39").Value = "'" + '5.376'
$ws.Range("E39").Value = "'" + '  +3.96%  '
$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D40").Value = "'" + '8.925'
$ws.Range("E40").Value = "'" + '  +2.79%  '
$ws.Range("D41").Value = "'" + '0.02343'
$ws.Range("E41").Value = "'" + '  +7.41%  '
$ws.Range("E42").Value = "'" + '  +2.41%  '
$ws.Range("D44").Value = "'" + '0.6183'
$ws.Range("E44").Value = "'" + '  +6.82%  '
$ws.Range("D45").Value = "'" + '0.9998'
$ws.Range("E45").Value = "'" + '  +1.33%  '
$ws.Range("D46").Value = "'" + '3.801'
$ws.Range("E46").Value = "'" + '  +4.57%  '
$ws.Range("D47").Value = "'" + '13.14'
$ws.Range("E47").Value = "'" + '  +1.68%  '
$ws.Range("D48").Value = "'" + '0.5976'
$ws.Range("E48").Value = "'" + '  +5.78%  '
$ws.Range("D49").Value = "'" + '126.88'
$ws.Range("E49").Value = "'" + '  +1.34%  '
$ws.Range("D50").Value = "'" + '2.035'
$ws.Range("E50").Value = "'" + '  +6.24%  '
$ws.Range("D51").Value = "'" + '0.07160'
$ws.Range("E51").Value = "'" + '  +5.68%  '
